$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "2018-10"
$ws.Cells.Item(2, 2).Value = 102.5
$ws.Cells.Item(2, 3).Value = 147.7
$ws.Cells.Item(3, 1).Value = "2018-11"
$ws.Cells.Item(3, 2).Value = 104
$ws.Cells.Item(3, 3).Value = 126.9
$ws.Cells.Item(4, 1).Value = "2018-12"
$ws.Cells.Item(4, 2).Value = 102.2
$ws.Cells.Item(4, 3).Value = 104.8
$ws.Cells.Item(5, 1).Value = "2018-01"
$ws.Cells.Item(5, 2).Value = 101.6564
$ws.Cells.Item(5, 3).Value = 113.5763
$ws.Cells.Item(6, 1).Value = "2018-02"
$ws.Cells.Item(6, 2).Value = 101.7
$ws.Cells.Item(6, 3).Value = 112.7
$ws.Cells.Item(7, 1).Value = "2018-03"
$ws.Cells.Item(7, 2).Value = 100.6
$ws.Cells.Item(7, 3).Value = 107.5
$ws.Cells.Item(8, 1).Value = "2018-04"
$ws.Cells.Item(8, 2).Value = 99
$ws.Cells.Item(8, 3).Value = 117
$ws.Cells.Item(9, 1).Value = "2018-05"
$ws.Cells.Item(9, 2).Value = 98.59999999999999
$ws.Cells.Item(9, 3).Value = 127.2
$ws.Cells.Item(10, 1).Value = "2018-06"
$ws.Cells.Item(10, 2).Value = 98.90000000000001
$ws.Cells.Item(10, 3).Value = 136.9
$ws.Cells.Item(11, 1).Value = "2018-07"
$ws.Cells.Item(11, 2).Value = 100.8
$ws.Cells.Item(11, 3).Value = 147.3
$ws.Cells.Item(12, 1).Value = "2018-08"
$ws.Cells.Item(12, 2).Value = 100.8
$ws.Cells.Item(12, 3).Value = 144.4
$ws.Cells.Item(13, 1).Value = "2018-09"
$ws.Cells.Item(13, 2).Value = 101.1
$ws.Cells.Item(13, 3).Value = 146.2
$ws.Cells.Item(14, 1).Value = "2019-10"
$ws.Cells.Item(14, 2).Value = 101.2
$ws.Cells.Item(14, 3).Value = 80.40000000000001
$ws.Cells.Item(15, 1).Value = "2019-11"
$ws.Cells.Item(15, 2).Value = 101.8
$ws.Cells.Item(15, 3).Value = 87.5
$ws.Cells.Item(16, 1).Value = "2019-12"
$ws.Cells.Item(16, 2).Value = 102.3
$ws.Cells.Item(16, 3).Value = 106.3
$ws.Cells.Item(17, 1).Value = "2019-01"
$ws.Cells.Item(17, 2).Value = 103.1
$ws.Cells.Item(17, 3).Value = 93.90000000000001
$ws.Cells.Item(18, 1).Value = "2019-02"
$ws.Cells.Item(18, 2).Value = 103.1
$ws.Cells.Item(18, 3).Value = 98.8
$ws.Cells.Item(19, 1).Value = "2019-03"
$ws.Cells.Item(19, 2).Value = 103.9
$ws.Cells.Item(19, 3).Value = 110.5
$ws.Cells.Item(20, 1).Value = "2019-04"
$ws.Cells.Item(20, 2).Value = 104.8
$ws.Cells.Item(20, 3).Value = 110.9
$ws.Cells.Item(21, 1).Value = "2019-05"
$ws.Cells.Item(21, 2).Value = 105.4
$ws.Cells.Item(21, 3).Value = 107.2
$ws.Cells.Item(22, 1).Value = "2019-06"
$ws.Cells.Item(22, 2).Value = 105.3
$ws.Cells.Item(22, 3).Value = 97.7
$ws.Cells.Item(23, 1).Value = "2019-07"
$ws.Cells.Item(23, 2).Value = 103.2
$ws.Cells.Item(23, 3).Value = 90.5
$ws.Cells.Item(24, 1).Value = "2019-08"
$ws.Cells.Item(24, 2).Value = 103.1
$ws.Cells.Item(24, 3).Value = 89.7
$ws.Cells.Item(25, 1).Value = "2019-09"
$ws.Cells.Item(25, 2).Value = 102.9
$ws.Cells.Item(25, 3).Value = 85
$ws.Cells.Item(26, 1).Value = "2020-10"
$ws.Cells.Item(26, 2).Value = 99.09999999999999
$ws.Cells.Item(26, 3).Value = 65.40000000000001
$ws.Cells.Item(27, 1).Value = "2020-11"
$ws.Cells.Item(27, 2).Value = 98.8
$ws.Cells.Item(27, 3).Value = 65.90000000000001
$ws.Cells.Item(28, 1).Value = "2020-12"
$ws.Cells.Item(28, 2).Value = 101.1
$ws.Cells.Item(28, 3).Value = 69
$ws.Cells.Item(29, 1).Value = "2020-01"
$ws.Cells.Item(29, 2).Value = 100.8
$ws.Cells.Item(29, 3).Value = 120.2
$ws.Cells.Item(30, 1).Value = "2020-02"
$ws.Cells.Item(30, 2).Value = 100.9
$ws.Cells.Item(30, 3).Value = 99.40000000000001
$ws.Cells.Item(31, 1).Value = "2020-03"
$ws.Cells.Item(31, 2).Value = 101
$ws.Cells.Item(31, 3).Value = 75
$ws.Cells.Item(32, 1).Value = "2020-04"
$ws.Cells.Item(32, 2).Value = 99.8
$ws.Cells.Item(32, 3).Value = 41.6
$ws.Cells.Item(33, 1).Value = "2020-05"
$ws.Cells.Item(33, 2).Value = 100
$ws.Cells.Item(33, 3).Value = 34.9
$ws.Cells.Item(34, 1).Value = "2020-06"
$ws.Cells.Item(34, 2).Value = 99.8
$ws.Cells.Item(34, 3).Value = 55.6
$ws.Cells.Item(35, 1).Value = "2020-07"
$ws.Cells.Item(35, 2).Value = 99.40000000000001
$ws.Cells.Item(35, 3).Value = 68.2
$ws.Cells.Item(36, 1).Value = "2020-08"
$ws.Cells.Item(36, 2).Value = 99.7
$ws.Cells.Item(36, 3).Value = 71.2
$ws.Cells.Item(37, 1).Value = "2020-09"
$ws.Cells.Item(37, 2).Value = 98.09999999999999
$ws.Cells.Item(37, 3).Value = 70.3
$ws.Cells.Item(38, 1).Value = "2021-10"
$ws.Cells.Item(38, 2).Value = 101.4
$ws.Cells.Item(38, 3).Value = 178.7
$ws.Cells.Item(39, 1).Value = "2021-11"
$ws.Cells.Item(39, 2).Value = 103.2
$ws.Cells.Item(39, 3).Value = 190
$ws.Cells.Item(40, 1).Value = "2021-12"
$ws.Cells.Item(40, 2).Value = 100.8
$ws.Cells.Item(40, 3).Value = 159.4
$ws.Cells.Item(41, 1).Value = "2021-01"
$ws.Cells.Item(41, 2).Value = 101.9
$ws.Cells.Item(41, 3).Value = 73.40000000000001
$ws.Cells.Item(42, 1).Value = "2021-02"
$ws.Cells.Item(42, 2).Value = 99.8
$ws.Cells.Item(42, 3).Value = 92.2
$ws.Cells.Item(43, 1).Value = "2021-03"
$ws.Cells.Item(43, 2).Value = 99.8
$ws.Cells.Item(43, 3).Value = 130.5
$ws.Cells.Item(44, 1).Value = "2021-04"
$ws.Cells.Item(44, 2).Value = 97.8
$ws.Cells.Item(44, 3).Value = 228.2
$ws.Cells.Item(45, 1).Value = "2021-05"
$ws.Cells.Item(45, 2).Value = 97.09999999999999
$ws.Cells.Item(45, 3).Value = 252.8
$ws.Cells.Item(46, 1).Value = "2021-06"
$ws.Cells.Item(46, 2).Value = 97.8
$ws.Cells.Item(46, 3).Value = 173.1
$ws.Cells.Item(47, 1).Value = "2021-07"
$ws.Cells.Item(47, 2).Value = 98.90000000000001
$ws.Cells.Item(47, 3).Value = 163.3
$ws.Cells.Item(48, 1).Value = "2021-08"
$ws.Cells.Item(48, 2).Value = 99.3
$ws.Cells.Item(48, 3).Value = 153.9
$ws.Cells.Item(49, 1).Value = "2021-09"
$ws.Cells.Item(49, 2).Value = 101.2
$ws.Cells.Item(49, 3).Value = 156.6
$ws.Cells.Item(50, 1).Value = "2022-10"
$ws.Cells.Item(50, 2).Value = 104.3
$ws.Cells.Item(50, 3).Value = 126.5
$ws.Cells.Item(51, 1).Value = "2022-11"
$ws.Cells.Item(51, 2).Value = 101.3
$ws.Cells.Item(51, 3).Value = 121.1
$ws.Cells.Item(52, 1).Value = "2022-12"
$ws.Cells.Item(52, 2).Value = 102.5
$ws.Cells.Item(52, 3).Value = 118.1
$ws.Cells.Item(53, 1).Value = "2022-01"
$ws.Cells.Item(53, 2).Value = 100.5
$ws.Cells.Item(53, 3).Value = 149.2
$ws.Cells.Item(54, 1).Value = "2022-02"
$ws.Cells.Item(54, 2).Value = 102.6
$ws.Cells.Item(54, 3).Value = 153.8
$ws.Cells.Item(55, 1).Value = "2022-03"
$ws.Cells.Item(55, 2).Value = 103.2
$ws.Cells.Item(55, 3).Value = 161.2
$ws.Cells.Item(56, 1).Value = "2022-04"
$ws.Cells.Item(56, 2).Value = 106.5
$ws.Cells.Item(56, 3).Value = 161.6
$ws.Cells.Item(57, 1).Value = "2022-05"
$ws.Cells.Item(57, 2).Value = 107
$ws.Cells.Item(57, 3).Value = 160.5
$ws.Cells.Item(58, 1).Value = "2022-06"
$ws.Cells.Item(58, 2).Value = 106.8
$ws.Cells.Item(58, 3).Value = 169.1
$ws.Cells.Item(59, 1).Value = "2022-07"
$ws.Cells.Item(59, 2).Value = 107.4
$ws.Cells.Item(59, 3).Value = 155.2
$ws.Cells.Item(60, 1).Value = "2022-08"
$ws.Cells.Item(60, 2).Value = 105.5
$ws.Cells.Item(60, 3).Value = 144.2
$ws.Cells.Item(61, 1).Value = "2022-09"
$ws.Cells.Item(61, 2).Value = 104.8
$ws.Cells.Item(61, 3).Value = 139.2
$ws.Cells.Item(62, 1).Value = "2023-01"
$ws.Cells.Item(62, 2).Value = 102.6
$ws.Cells.Item(62, 3).Value = 105.4
$ws.Cells.Item(63, 1).Value = "2023-02"
$ws.Cells.Item(63, 2).Value = 103.9
$ws.Cells.Item(63, 3).Value = 94.59999999999999
$ws.Cells.Item(64, 1).Value = "2023-03"
$ws.Cells.Item(64, 2).Value = 103
$ws.Cells.Item(64, 3).Value = 79.5
$ws.Cells.Item(65, 1).Value = "2023-04"
$ws.Cells.Item(65, 2).Value = 103.2
$ws.Cells.Item(65, 3).Value = 79
$ws.Cells.Item(66, 1).Value = "2023-05"
$ws.Cells.Item(66, 2).Value = 102.7
$ws.Cells.Item(66, 3).Value = 75.59999999999999
$ws.Cells.Item(67, 1).Value = "2023-06"
$ws.Cells.Item(67, 2).Value = 102.1
$ws.Cells.Item(67, 3).Value = 68.2
$ws.Cells.Item(68, 1).Value = "2023-07"
$ws.Cells.Item(68, 2).Value = 101.3
$ws.Cells.Item(68, 3).Value = 73.5

# Ensure column A keeps the original centered/bordered style (s=1) for all data rows,
# including the newly appended ones.
$ws.Range("A2").Copy()
$ws.Range("A2:A68").PasteSpecial(-4122)
$excel.CutCopyMode = $false
